$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, shifting existing rows 104-163 down to 105-164
$ws.Rows("104:104").Insert()

# The new row 104 should duplicate the (now shifted) row 105 for the columns that
# stay constant, and carry the new values for Fecha/Volumen/Precios/Precio $/Kg.
$ws.Cells.Item(104, 1).Value = 11
$ws.Cells.Item(104, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(104, 3).Value = "Bíobío"
$ws.Cells.Item(104, 4).Value = 44873
$ws.Cells.Item(104, 5).Value = 8
$ws.Cells.Item(104, 6).Value = 100112043
$ws.Cells.Item(104, 7).Value = "Pepino ensalada"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 220
$ws.Cells.Item(104, 11).Value = 15000
$ws.Cells.Item(104, 12).Value = 17000
$ws.Cells.Item(104, 13).Value = 15909
$ws.Cells.Item(104, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(104, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(104, 16).Value = 265
$ws.Cells.Item(104, 17).Value = 60
$ws.Cells.Item(104, 18).Value = "Hortaliza"
